$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: shorten "Hillary Clinton" to "Hillary"
$ws.Range("B1").Value = "Hillary"

# New section: Favorite Tech Device (B17 written before J1 to match original authoring order)
$ws.Range("B17").Value = "blackberry"

# Row 1: add new candidate "lindsay graham"
$ws.Range("J1").Value = "lindsay graham"

$ws.Range("J17").Value = "iphone"
$ws.Range("A17").Value = "Favorite Tech Device"

# New section: favorite drinks
$ws.Range("A18").Value = "favorite drinks"
$ws.Range("B18").Value = "wine"
$ws.Range("G18").Value = "water"

# Update active cell selection
[void]$ws.Range("E6").Select()
